# "fix(gui) step 1 and 2"
# Step 1: roll the price-list date in A1 forward by one day.
# Step 2: update the three fratacho prices (D29:D31) to the new prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1 - bump the date serial in A1 (17-Jan-2024 -> 18-Jan-2024)
$ws.Range("A1").Value = 45309

# Step 2 - new prices for FRATACHO Nº 20 / Nº 25 / Nº 30
$ws.Range("D29").Value = 598
$ws.Range("D30").Value = 640
$ws.Range("D31").Value = 815
